# "excel actualizados a GKD only" — drop the non-GKD benchmark instances
# (MDG-a, MDG-b, SOM-a x2, SOM-b) from the results table, leaving only the
# GKD-* rows. The two summary rows (5 and 6) use AVERAGE(...9:...32), so once
# the MDG/SOM rows are blanked those averages automatically recompute over
# just the remaining GKD data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 28:32 held MDG-a_9, MDG-b_12, SOM-a_18, SOM-a_33 and SOM-b_2 — remove
# the instance name (col A) entirely and clear out all of their metric
# values (cols B:O), leaving just the blank, still-styled cells behind.
$ws.Range("A28:A32").ClearContents()
$ws.Range("B28:O32").ClearContents()

# Reflect the author's final selection position on the sheet.
$ws.Range("P25").Select()
